# Weekly fruit/vegetable price update: insert two new price records
# (dated 2023-03-30 / serial 45015) at the top of the existing block for
# "Vega Modelo de Temuco - Coliflor" and push the rest of the rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 506; everything below (old rows
# 506-576) shifts down to rows 508-578, carrying its data/formatting
# along with it.
$ws.Rows("506:507").Insert()

# --- New row 506 ---
$ws.Range("A506").Value = 10
$ws.Range("B506").Value = "Vega Modelo de Temuco"
$ws.Range("C506").Value = "La Araucanía"
$ws.Range("D506").Value = 45015
$ws.Range("E506").Value = 9
$ws.Range("F506").Value = 100112008
$ws.Range("G506").Value = "Coliflor"
$ws.Range("H506").Value = "Sin especificar"
$ws.Range("I506").Value = "Primera"
$ws.Range("J506").Value = 500
$ws.Range("K506").Value = 1300
$ws.Range("L506").Value = 1400
$ws.Range("M506").Value = 1340
$ws.Range("N506").Value = "$/unidad"
$ws.Range("O506").Value = "Provincia de Cautín"
$ws.Range("P506").Value = 1340
$ws.Range("Q506").Value = 1
$ws.Range("R506").Value = "Hortaliza"

# --- New row 507 ---
$ws.Range("A507").Value = 10
$ws.Range("B507").Value = "Vega Modelo de Temuco"
$ws.Range("C507").Value = "La Araucanía"
$ws.Range("D507").Value = 45015
$ws.Range("E507").Value = 9
$ws.Range("F507").Value = 100112008
$ws.Range("G507").Value = "Coliflor"
$ws.Range("H507").Value = "Sin especificar"
$ws.Range("I507").Value = "Primera"
$ws.Range("J507").Value = 2500
$ws.Range("K507").Value = 1300
$ws.Range("L507").Value = 1300
$ws.Range("M507").Value = 1300
$ws.Range("N507").Value = "$/unidad"
$ws.Range("O507").Value = "Región Metropolitana"
$ws.Range("P507").Value = 1300
$ws.Range("Q507").Value = 1
$ws.Range("R507").Value = "Hortaliza"
